$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "capacity_to_activity" parameter row is inserted for the biomass CHP
# entity right after the existing "buildrate" row (old row 9). Inserting a
# whole row shifts every row from the old row 10 downwards (old rows 10-316)
# down by one (new rows 11-317), matching the rest of the table.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row with the new entity-dependent
# capacity-to-activity conversion factor.
$ws.Range("A10").Value2 = "CHE"
$ws.Range("B10").Value2 = "conv_chp_biomass"
$ws.Range("C10").Value2 = "capacity_to_activity"
$ws.Range("D10").Value2 = "constant"
$ws.Range("G10").Value2 = 0.001
$ws.Range("H10").Value2 = "GW/TWh"

# The used range grew by one row (table now spans down to row 850 instead of
# 849), so refresh the autofilter to cover the extra row.
$ws.AutoFilterMode = $false
$ws.Range("A5:L850").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the autofilter.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$850"
    }
}

# Match the author's final cursor position.
$ws.Activate()
$ws.Range("C12").Select()
